$d = $word.ActiveDocument

# --- 1) Strike through "Put a blanking screen over each box when not
#        clicked or hovered over" (paragraph text + paragraph mark). ---
$rng1 = $d.Content
$found1 = $rng1.Find.Execute(
    "Put a blanking screen over each box when not clicked or hovered over",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $para1 = $rng1.Paragraphs(1)
    $para1.Range.Font.StrikeThrough = 1
}

# --- 2) Strike through "When question mark box is clicked, image stays
#        revealed" (paragraph text + paragraph mark). ---
$rng2 = $d.Content
$found2 = $rng2.Find.Execute(
    "When question mark box is clicked, image stays revealed",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $para2 = $rng2.Paragraphs(1)
    $para2.Range.Font.StrikeThrough = 1
}

# --- 3) Split the run "Create an input field where the user puts in a
#        name that becomes a parameter" into two runs, with the
#        _GoBack bookmark relocated between them (right after "th"). ---
$rng3 = $d.Content
$found3 = $rng3.Find.Execute(
    "Create an input field where th",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found3) {
    $splitPoint = $rng3.End
    $splitRange = $d.Range($splitPoint, $splitPoint)
    $d.Bookmarks.Add("_GoBack", $splitRange)
}
